# Tutorial 6 solution update: change the date separator in column A
# (rows 3-21) from "/" to "-", e.g. "28/07/2022" -> "28-07-2022".
#
# The dates are stored as plain text. Writing an ambiguous day/month
# string (day <= 12) straight into Range.Value2 makes Excel "smart"
# re-type it as a real date serial (and mint a date number format),
# which is not what the source data wants. Routing the new text
# through a formula and then collapsing it to a literal via
# Copy + PasteSpecial(xlPasteValues) keeps it as plain text and
# leaves the cell's existing style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    $cell.Formula = '="' + $dates[$row] + '"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}
